$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (GitHub Actions scheduled update).
# Column D ("Price") cells that would parse as plain numbers must be
# pre-formatted as Text so Excel keeps storing them as text, matching the
# source data feed (which always writes Price/Volume as text).
$numericPriceRows = @(5, 9, 10, 11, 14, 16, 19, 22, 24, 25, 26, 27, 28, 29, 33, 41, 44, 46, 49, 50)
foreach ($r in $numericPriceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Column D ("Price") updates
$ws.Cells.Item(2, 4).Value = "26.801.13"
$ws.Cells.Item(3, 4).Value = "1.638.28"
$ws.Cells.Item(5, 4).Value = "218.62"
$ws.Cells.Item(9, 4).Value = "0.0622"
$ws.Cells.Item(10, 4).Value = "19.23"
$ws.Cells.Item(11, 4).Value = "0.0844"
$ws.Cells.Item(12, 4).Value = "1.866.72"
$ws.Cells.Item(13, 4).Value = "1.657.15"
$ws.Cells.Item(14, 4).Value = "4.14"
$ws.Cells.Item(16, 4).Value = "64.75"
$ws.Cells.Item(17, 4).Value = "26.792.76"
$ws.Cells.Item(19, 4).Value = "214.60"
$ws.Cells.Item(22, 4).Value = "6.44"
$ws.Cells.Item(24, 4).Value = "9.14"
$ws.Cells.Item(25, 4).Value = "147.59"
$ws.Cells.Item(26, 4).Value = "1.00"
$ws.Cells.Item(27, 4).Value = "0.119"
$ws.Cells.Item(28, 4).Value = "7.03"
$ws.Cells.Item(29, 4).Value = "15.71"
$ws.Cells.Item(33, 4).Value = "2.98"
$ws.Cells.Item(35, 4).Value = "1.261.67"
$ws.Cells.Item(41, 4).Value = "0.805"
$ws.Cells.Item(43, 4).Value = "1.778.10"
$ws.Cells.Item(44, 4).Value = "2.14"
$ws.Cells.Item(46, 4).Value = "60.10"
$ws.Cells.Item(49, 4).Value = "0.0963"
$ws.Cells.Item(50, 4).Value = "7.54"

# Column E ("Volume(1h)") updates
$ws.Cells.Item(2, 5).Value = "  +0.07%  "
$ws.Cells.Item(3, 5).Value = "  -0.36%  "
$ws.Cells.Item(4, 5).Value = "  -0.26%  "
$ws.Cells.Item(6, 5).Value = "  -0.45%  "
$ws.Cells.Item(7, 5).Value = "  -0.33%  "
$ws.Cells.Item(8, 5).Value = "  -0.59%  "
$ws.Cells.Item(9, 5).Value = "  -0.66%  "
$ws.Cells.Item(10, 5).Value = "  +0.05%  "
$ws.Cells.Item(11, 5).Value = "  +0.33%  "
$ws.Cells.Item(12, 5).Value = "  -0.36%  "
$ws.Cells.Item(13, 5).Value = "  +1.25%  "
$ws.Cells.Item(14, 5).Value = "  -1.01%  "
$ws.Cells.Item(15, 5).Value = "  -0.45%  "
$ws.Cells.Item(16, 5).Value = "  +0.20%  "
$ws.Cells.Item(17, 5).Value = "  +0.03%  "
$ws.Cells.Item(18, 5).Value = "  -0.91%  "
$ws.Cells.Item(19, 5).Value = "  +0.17%  "
$ws.Cells.Item(20, 5).Value = "  -0.26%  "
$ws.Cells.Item(21, 5).Value = "  -0.14%  "
$ws.Cells.Item(22, 5).Value = "  +2.40%  "
$ws.Cells.Item(23, 5).Value = "  -2.92%  "
$ws.Cells.Item(24, 5).Value = "  -2.51%  "
$ws.Cells.Item(25, 5).Value = "  +1.87%  "
$ws.Cells.Item(26, 5).Value = "  -0.44%  "
$ws.Cells.Item(27, 5).Value = "  +0.22%  "
$ws.Cells.Item(28, 5).Value = "  -0.99%  "
$ws.Cells.Item(29, 5).Value = "  +0.07%  "
$ws.Cells.Item(30, 5).Value = "  -1.82%  "
$ws.Cells.Item(31, 5).Value = "  +1.28%  "
$ws.Cells.Item(32, 5).Value = "  +1.51%  "
$ws.Cells.Item(33, 5).Value = "  -0.31%  "
$ws.Cells.Item(34, 5).Value = "  +0.15%  "
$ws.Cells.Item(35, 5).Value = "  -2.06%  "
$ws.Cells.Item(36, 5).Value = "  +0.24%  "
$ws.Cells.Item(37, 5).Value = "  -0.39%  "
$ws.Cells.Item(38, 5).Value = "  -2.44%  "
$ws.Cells.Item(39, 5).Value = "  -1.43%  "
$ws.Cells.Item(40, 5).Value = "  -0.26%  "
$ws.Cells.Item(41, 5).Value = "  -0.68%  "
$ws.Cells.Item(42, 5).Value = "  -0.33%  "
$ws.Cells.Item(43, 5).Value = "  -0.92%  "
$ws.Cells.Item(44, 5).Value = "  -4.22%  "
$ws.Cells.Item(45, 5).Value = "  +0.69%  "
$ws.Cells.Item(46, 5).Value = "  +0.26%  "
$ws.Cells.Item(47, 5).Value = "  -1.79%  "
$ws.Cells.Item(48, 5).Value = "  -1.04%  "
$ws.Cells.Item(49, 5).Value = "  -1.39%  "
$ws.Cells.Item(50, 5).Value = "  -1.76%  "
$ws.Cells.Item(51, 5).Value = "  -0.25%  "
